$d = $word.ActiveDocument

# Locate the paragraph that still holds the un-split "{m:userdoc 'zone1'}"
# field-ish text (originally stored as two runs: "{m" and ":userdoc 'zone1'}").
$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*{m*" -and $t -like "*:userdoc*zone1*") {
        $target = $p
    }
}

if ($target -ne $null) {
    $pRange = $target.Range
    # Exclude the trailing paragraph mark from the range we rewrite.
    $textRange = $d.Range($pRange.Start, $pRange.End - 1)

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' +
        '<w:p>' +
        '<w:r><w:t>{</w:t></w:r>' +
        '<w:r><w:t>m</w:t></w:r>' +
        '<w:r><w:t>:userdoc ''zone1''</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve">}</w:t></w:r>' +
        '</w:p>' +
        '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $textRange.InsertXML($xml)
}
